$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update 2017 fitness estimates: total fitness (W) should be sum of 2016 and 2017 seed mass
$ws.Range("C8").Value = 544.65099999999995
$ws.Range("D8").Value = 163.02699999999999

$ws.Range("D9").Value = 81.897999999999996

$ws.Range("C10").Value = 70.379000000000005
$ws.Range("D10").Value = 75.872
